$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 78.416664
$ws.Range("I9").Value = 64.75
$ws.Range("J9").Value = 105.75
$ws.Range("K9").Value = 64.75
$ws.Range("L9").Value = 105.75
$ws.Range("M9").Value = 104.25
$ws.Range("N9").Value = -443.75

# Row 99
$ws.Range("H99").Value = 1348.875

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 367.36365
$ws.Range("I5").Value = 284.57144
$ws.Range("J5").Value = 512.25
$ws.Range("K5").Value = 284.57144
$ws.Range("L5").Value = 512.25
$ws.Range("M5").Value = -172.57144
$ws.Range("N5").Value = -736.25

# Row 32
$ws.Range("H32").Value = 7723.242
$ws.Range("I32").Value = 6402.0938
$ws.Range("K32").Value = 6402.0938
$ws.Range("M32").Value = -6115.0938

# Row 61
$ws.Range("H61").Value = 5000
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5424

# Row 63
$ws.Range("H63").Value = 2115.4285
$ws.Range("J63").Value = 4000
$ws.Range("L63").Value = 4000
$ws.Range("N63").Value = -5372

# Row 66
$ws.Range("H66").Value = 2115.4285
$ws.Range("J66").Value = 4000
$ws.Range("L66").Value = 20000
$ws.Range("N66").Value = -26864

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# Row 122
$ws.Range("H122").Value = 2647.8333
$ws.Range("I122").Value = 2597.6
$ws.Range("K122").Value = 7792.799999999999
$ws.Range("M122").Value = -5342.799999999999

# Row 136
$ws.Range("H136").Value = 5000
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 367.36365
$ws.Range("I4").Value = 284.57144
$ws.Range("J4").Value = 512.25
$ws.Range("K4").Value = 284.57144
$ws.Range("L4").Value = 512.25
$ws.Range("M4").Value = -169.57144
$ws.Range("N4").Value = -742.25

# Row 107
$ws.Range("H107").Value = 1017.5
$ws.Range("I107").Value = 1017.5
$ws.Range("K107").Value = 1017.5
$ws.Range("M107").Value = 902.5

$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 19341
$ws.Range("J28").Value = 19341
$ws.Range("L28").Value = 19341
$ws.Range("N28").Value = -19831

# Row 58
$ws.Range("H58").Value = 4296.1
$ws.Range("I58").Value = 3247.5
$ws.Range("J58").Value = 4995.1665
$ws.Range("K58").Value = 3247.5
$ws.Range("L58").Value = 4995.1665
$ws.Range("M58").Value = -3044.5
$ws.Range("N58").Value = -5401.1665

# Row 59
$ws.Range("H59").Value = 30832.555
$ws.Range("I59").Value = 17501.5
$ws.Range("J59").Value = 34641.43
$ws.Range("K59").Value = 17501.5
$ws.Range("L59").Value = 34641.43
$ws.Range("M59").Value = -16356.5
$ws.Range("N59").Value = -36931.43

# Row 60
$ws.Range("H60").Value = 18782.285
$ws.Range("J60").Value = 24998
$ws.Range("L60").Value = 24998
$ws.Range("N60").Value = -26020

# Row 68
$ws.Range("H68").Value = 37804.668
$ws.Range("I68").Value = 20268
$ws.Range("K68").Value = 20268
$ws.Range("M68").Value = -19519

# Row 71
$ws.Range("H71").Value = 37804.668
$ws.Range("I71").Value = 20268
$ws.Range("K71").Value = 60804
$ws.Range("M71").Value = -57060

# Row 103
$ws.Range("H103").Value = 6662.25
$ws.Range("I103").Value = 6662.25
$ws.Range("K103").Value = 6662.25
$ws.Range("M103").Value = -5490.25

# Row 108
$ws.Range("H108").Value = 30000
$ws.Range("I108").Value = 30000
$ws.Range("K108").Value = 30000
$ws.Range("M108").Value = -26160

# Row 132
$ws.Range("H132").Value = 1381.6428
$ws.Range("I132").Value = 1381.6428
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4144.928400000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1614.928400000001
$ws.Range("N132").ClearContents()

# Row 136
$ws.Range("H136").Value = 4296.1
$ws.Range("I136").Value = 3247.5
$ws.Range("J136").Value = 4995.1665
$ws.Range("K136").Value = 9742.5
$ws.Range("L136").Value = 14985.4995
$ws.Range("M136").Value = -7192.5
$ws.Range("N136").Value = -20085.4995

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 97920
$ws.Range("J37").Value = 97920
$ws.Range("L37").Value = 293760
$ws.Range("N37").Value = -293984

# Row 113
$ws.Range("H113").Value = 975.3333
$ws.Range("I113").Value = 479.5
$ws.Range("J113").Value = 1223.25
$ws.Range("K113").Value = 1438.5
$ws.Range("L113").Value = 3669.75
$ws.Range("M113").Value = 731.5
$ws.Range("N113").Value = -8009.75

# Row 129
$ws.Range("H129").Value = 1339.4615
$ws.Range("J129").Value = 1820.5714
$ws.Range("L129").Value = 5461.7142
$ws.Range("N129").Value = -15461.7142

$ws = $wb.Worksheets.Item("GSM")
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Row 102
$ws.Range("H102").Value = 1447.5
$ws.Range("I102").Value = 1137
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1137
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 485
$ws.Range("N102").Value = -6244

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2276.6
$ws.Range("I16").Value = 2395.75
$ws.Range("K16").Value = 2395.75
$ws.Range("M16").Value = -2225.75

# Row 55
$ws.Range("H55").Value = 703.1
$ws.Range("I55").Value = 1042.2
$ws.Range("J55").Value = 364
$ws.Range("K55").Value = 1042.2
$ws.Range("L55").Value = 364
$ws.Range("M55").Value = -869.2
$ws.Range("N55").Value = -710

# Row 61
$ws.Range("H61").Value = 774.75
$ws.Range("I61").Value = 774.75
$ws.Range("K61").Value = 774.75
$ws.Range("M61").Value = -572.75

# Row 93
$ws.Range("H93").Value = 147.5
$ws.Range("I93").Value = 147.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 147.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 1100.5
$ws.Range("N93").ClearContents()

# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# Row 113
$ws.Range("H113").Value = 774.75
$ws.Range("I113").Value = 774.75
$ws.Range("K113").Value = 774.75
$ws.Range("M113").Value = 1395.25

# Row 132
$ws.Range("H132").Value = 7120.421
$ws.Range("I132").Value = 6959.6665
$ws.Range("K132").Value = 20878.9995
$ws.Range("M132").Value = -18348.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4997
$ws.Range("I126").Value = 4997
$ws.Range("K126").Value = 14991
$ws.Range("M126").Value = -12521
